$d = $word.ActiveDocument

$replacements = @(
    @{ old = "Найти имена";    new = "Щицбц фхсьн" },
    @{ old = "Найди отличия";  new = "Щицтц ъышчефз" },
    @{ old = "Найди предметы"; new = "Щицтц ыщстърыз" },
    @{ old = "Найти букву";    new = "Щицбц мьчрб" },
    @{ old = "Лабиринт";       new = "Чинчюфця" },
    @{ old = "Математика";     new = "Шияуълыхщн" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
